$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 169 (new weekly data for
# Vega Monumental Concepción - Acelga), pushing the existing rows 169-174
# down to 171-176.
$ws.Range("A169:R170").EntireRow.Insert()

# Row 169: "Primera" quality, new date, Región de Ñuble
$ws.Cells.Item(169, 1).Value = 11
$ws.Cells.Item(169, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(169, 3).Value = "Bíobío"
$ws.Cells.Item(169, 4).Value = 44509
$ws.Cells.Item(169, 5).Value = 8
$ws.Cells.Item(169, 6).Value = 100112009
$ws.Cells.Item(169, 7).Value = "Acelga"
$ws.Cells.Item(169, 8).Value = "Sin especificar"
$ws.Cells.Item(169, 9).Value = "Primera"
$ws.Cells.Item(169, 10).Value = 200
$ws.Cells.Item(169, 11).Value = 600
$ws.Cells.Item(169, 12).Value = 700
$ws.Cells.Item(169, 13).Value = 650
$ws.Cells.Item(169, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(169, 15).Value = "Región de Ñuble"
$ws.Cells.Item(169, 16).Value = 650
$ws.Cells.Item(169, 17).Value = 1
$ws.Cells.Item(169, 18).Value = "Hortaliza"

# Row 170: "Segunda" quality, new date, Región de Ñuble
$ws.Cells.Item(170, 1).Value = 11
$ws.Cells.Item(170, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(170, 3).Value = "Bíobío"
$ws.Cells.Item(170, 4).Value = 44509
$ws.Cells.Item(170, 5).Value = 8
$ws.Cells.Item(170, 6).Value = 100112009
$ws.Cells.Item(170, 7).Value = "Acelga"
$ws.Cells.Item(170, 8).Value = "Sin especificar"
$ws.Cells.Item(170, 9).Value = "Segunda"
$ws.Cells.Item(170, 10).Value = 100
$ws.Cells.Item(170, 11).Value = 500
$ws.Cells.Item(170, 12).Value = 500
$ws.Cells.Item(170, 13).Value = 500
$ws.Cells.Item(170, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(170, 15).Value = "Región de Ñuble"
$ws.Cells.Item(170, 16).Value = 500
$ws.Cells.Item(170, 17).Value = 1
$ws.Cells.Item(170, 18).Value = "Hortaliza"
